$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,10
$arr[0,0] = 8.950324975496269
$arr[0,1] = 7.363184871255143
$arr[0,2] = 13.19849402262084
$arr[0,3] = 42.68664402002148
$arr[0,4] = 50.9240849741746
$arr[0,5] = 19.92095029115883
$arr[0,6] = 32.97184719964956
$arr[0,7] = 10.6380506989057
$arr[0,8] = 21.09292738975575
$arr[0,9] = 10.46520630913917
$arr[1,0] = 8.890693630935456
$arr[1,1] = 7.315966163578727
$arr[1,2] = 13.19054907684496
$arr[1,3] = 42.86244418527542
$arr[1,4] = 51.18971017092255
$arr[1,5] = 20.01788411387478
$arr[1,6] = 33.11653979976615
$arr[1,7] = 10.66209144419277
$arr[1,8] = 20.60778879199564
$arr[1,9] = 10.48220430971549
$arr[2,0] = 8.854868360279134
$arr[2,1] = 7.287213331911786
$arr[2,2] = 13.18735302576636
$arr[2,3] = 42.98247200471939
$arr[2,4] = 51.37086215606471
$arr[2,5] = 20.08177748253189
$arr[2,6] = 33.21403736568448
$arr[2,7] = 10.67805437236292
$arr[2,8] = 20.30615757025016
$arr[2,9] = 10.49371978272552
$arr[3,0] = 8.840477907332895
$arr[3,1] = 7.275563467309532
$arr[3,2] = 13.18647526498306
$arr[3,3] = 43.03440979836601
$arr[3,4] = 51.44919388666914
$arr[3,5] = 20.10891304141727
$arr[3,6] = 33.25593687791794
$arr[3,7] = 10.68486193606472
$arr[3,8] = 20.18245780067639
$arr[3,9] = 10.49868402568438
$arr[4,0] = 8.838101280495888
$arr[4,1] = 7.273633270699663
$arr[4,2] = 13.18635520705063
$arr[4,3] = 43.0432163582791
$arr[4,4] = 51.46247237613303
$arr[4,5] = 20.11348517506086
$arr[4,6] = 33.26302499177742
$arr[4,7] = 10.68601060939822
$arr[4,8] = 20.16187513263587
$arr[4,9] = 10.49952474859643
$arr[5,0] = 8.854673427404203
$arr[5,1] = 7.287055936856039
$arr[5,2] = 13.18733946647991
$arr[5,3] = 42.98316022419984
$arr[5,4] = 51.37190034081589
$arr[5,5] = 20.08213899680728
$arr[5,6] = 33.21459366797818
$arr[5,7] = 10.67814495613149
$arr[5,8] = 20.30449226230111
$arr[5,9] = 10.4937856320924
$arr[6,0] = 8.92960669700466
$arr[6,1] = 7.346858283827083
$arr[6,2] = 13.19540650695394
$arr[6,3] = 42.74474342146723
$arr[6,4] = 51.011907026074
$arr[6,5] = 19.95346383095642
$arr[6,6] = 33.01993587447267
$arr[6,7] = 10.64609071330558
$arr[6,8] = 20.9265254193253
$arr[6,9] = 10.4708435695032
$arr[7,0] = 9.082342116239978
$arr[7,1] = 7.465726384800497
$arr[7,2] = 13.22449567287344
$arr[7,3] = 42.373699799849
$arr[7,4] = 50.45054740481847
$arr[7,5] = 19.73593168938777
$arr[7,6] = 32.70725021773899
$arr[7,7] = 10.59275345362597
$arr[7,8] = 22.10982858381671
$arr[7,9] = 10.43439675676555
$arr[8,0] = 9.197449101120347
$arr[8,1] = 7.553619130013114
$arr[8,2] = 13.25384602847249
$arr[8,3] = 42.1607132428422
$arr[8,4] = 50.12802309977297
$arr[8,5] = 19.59743792663054
$arr[8,6] = 32.52009250024534
$arr[8,7] = 10.55935071489566
$arr[8,8] = 22.94854048610991
$arr[8,9] = 10.41280556783955
$arr[9,0] = 9.250302691051647
$arr[9,1] = 7.593636709562119
$arr[9,2] = 13.26890386946911
$arr[9,3] = 42.07692908735626
$arr[9,4] = 50.00118697870541
$arr[9,5] = 19.53908979591331
$arr[9,6] = 32.44429585022823
$arr[9,7] = 10.54540667717452
$arr[9,8] = 23.32183333815084
$arr[9,9] = 10.40410466247471
$arr[10,0] = 9.270375584727136
$arr[10,1] = 7.608787758696938
$arr[10,2] = 13.27484864184345
$arr[10,3] = 42.04709943045726
$arr[10,4] = 49.95604557336534
$arr[10,5] = 19.51766631421065
$arr[10,6] = 32.41694532338985
$arr[10,7] = 10.54030600388497
$arr[10,8] = 23.46188353038645
$arr[10,9] = 10.40097066338296
$arr[11,0] = 9.266050125888007
$arr[11,1] = 7.605524954009319
$arr[11,2] = 13.27355757887962
$arr[11,3] = 42.05343918044743
$arr[11,4] = 49.96563864347483
$arr[11,5] = 19.52225032969691
$arr[11,6] = 32.42277548186095
$arr[11,7] = 10.54139654023264
$arr[11,8] = 23.43178122829421
$arr[11,9] = 10.40163847834406
$arr[12,0] = 9.251952999979597
$arr[12,1] = 7.594883294645123
$arr[12,2] = 13.26938810130467
$arr[12,3] = 42.07443687933583
$arr[12,4] = 49.99741510740375
$arr[12,5] = 19.53731379678437
$arr[12,6] = 32.44201856069701
$arr[12,7] = 10.54498344316157
$arr[12,8] = 23.33338222893365
$arr[12,9] = 10.40384360499652
$arr[13,0] = 9.24332534408728
$arr[13,1] = 7.588364384931941
$arr[13,2] = 13.26686570027653
$arr[13,3] = 42.08754608065419
$arr[13,4] = 50.01725616838214
$arr[13,5] = 19.54662815540674
$arr[13,6] = 32.45398182775093
$arr[13,7] = 10.5472039101573
$arr[13,8] = 23.27293622057312
$arr[13,9] = 10.4052152447091
$arr[14,0] = 9.19400383048333
$arr[14,1] = 7.551003924574455
$arr[14,2] = 13.2528960732037
$arr[14,3] = 42.16645356715745
$arr[14,4] = 50.13671476836047
$arr[14,5] = 19.60134496776368
$arr[14,6] = 32.52523471999516
$arr[14,7] = 10.56028714182606
$arr[14,8] = 22.92396838594732
$arr[14,9] = 10.41339672011107
$arr[15,0] = 9.163863738703567
$arr[15,1] = 7.528087882199324
$arr[15,2] = 13.24476132417622
$arr[15,3] = 42.21822664953196
$arr[15,4] = 50.21511233778619
$arr[15,5] = 19.63610558946461
$arr[15,6] = 32.57134502436171
$arr[15,7] = 10.56863351139955
$arr[15,8] = 22.70768686850151
$arr[15,9] = 10.41870266188837
$arr[16,0] = 9.146574890734152
$arr[16,1] = 7.514910912357699
$arr[16,2] = 13.24024319840957
$arr[16,3] = 42.2492378553249
$arr[16,4] = 50.2620741870242
$arr[16,5] = 19.65653682217996
$arr[16,6] = 32.59874529112998
$arr[16,7] = 10.57355189175532
$arr[16,8] = 22.58251847043143
$arr[16,9] = 10.42186003838081
$arr[17,0] = 9.140729615758435
$arr[17,1] = 7.510450293486827
$arr[17,2] = 13.2387411287587
$arr[17,3] = 42.25994902473889
$arr[17,4] = 50.27829485676888
$arr[17,5] = 19.66352960990475
$arr[17,6] = 32.60817322752029
$arr[17,7] = 10.57523740593371
$arr[17,8] = 22.54001041741651
$arr[17,9] = 10.42294720902319
$arr[18,0] = 9.167067433415745
$arr[18,1] = 7.530527003851782
$arr[18,2] = 13.24561066208624
$arr[18,3] = 42.21258764039029
$arr[18,4] = 50.20657306478829
$arr[18,5] = 19.63235992885249
$arr[18,6] = 32.56634548468837
$arr[18,7] = 10.56773283938577
$arr[18,8] = 22.73079091090257
$arr[18,9] = 10.41812691501602
$arr[19,0] = 9.2560921804352
$arr[19,1] = 7.598009143874449
$arr[19,2] = 13.27060621325715
$arr[19,3] = 42.06821774362523
$arr[19,4] = 49.98800295898984
$arr[19,5] = 19.53287104138247
$arr[19,6] = 32.43632964073408
$arr[19,7] = 10.5439250098057
$arr[19,8] = 23.3623208214566
$arr[19,9] = 10.40319154346159
$arr[20,0] = 9.31460963659106
$arr[20,1] = 7.642093081897242
$arr[20,2] = 13.28835551466155
$arr[20,3] = 41.98493001674105
$arr[20,4] = 49.86200539420629
$arr[20,5] = 19.47176563925005
$arr[20,6] = 32.35924125811619
$arr[20,7] = 10.52941212141186
$arr[20,8] = 23.7673895548705
$arr[20,9] = 10.3943677694462
$arr[21,0] = 9.283351189247115
$arr[21,1] = 7.618568992568966
$arr[21,2] = 13.27875398457842
$arr[21,3] = 42.02836544380517
$arr[21,4] = 49.92770123194922
$arr[21,5] = 19.50401954347877
$arr[21,6] = 32.3996605446802
$arr[21,7] = 10.5370622191681
$arr[21,8] = 23.55193686215036
$arr[21,9] = 10.39899153645277
$arr[22,0] = 9.165618920361368
$arr[22,1] = 7.529424283332372
$arr[22,2] = 13.24522618205725
$arr[22,3] = 42.21513315571563
$arr[22,4] = 50.21042778506622
$arr[22,5] = 19.63405195017388
$arr[22,6] = 32.5686030023413
$arr[22,7] = 10.56813965958107
$arr[22,8] = 22.72034813445874
$arr[22,9] = 10.41838687711063
$arr[23,0] = 9.040468202028753
$arr[23,1] = 7.433440893332283
$arr[23,2] = 13.21521647267866
$arr[23,3] = 42.46366334974645
$arr[23,4] = 50.58674697327632
$arr[23,5] = 19.79104297986325
$arr[23,6] = 32.78439851582147
$arr[23,7] = 10.6061653878717
$arr[23,8] = 21.79452374243852
$arr[23,9] = 10.44334423807297

$ws.Range("C2:L25").Value = $arr
Write-Output "done"
